$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the C column values for rows 2-15 (new sampled data, test structure before changes)
$ws.Range("C2").Value = 0.3285138811888638
$ws.Range("C3").Value = 0.3216747188243654
$ws.Range("C4").Value = 0.3286024557405686
$ws.Range("C5").Value = 0.3211807332839862
$ws.Range("C6").Value = 0.3288331149508074
$ws.Range("C7").Value = 0.3219496020622352
$ws.Range("C8").Value = 0.3285318265659036
$ws.Range("C9").Value = 0.3224131462840029
$ws.Range("C10").Value = 0.3287359510336123
$ws.Range("C11").Value = 0.3230641013449289
$ws.Range("C12").Value = 0.3286534576049262
$ws.Range("C13").Value = 0.3240327853785042
$ws.Range("C14").Value = 0.3287990759682234
$ws.Range("C15").Value = 0.3201716152902284

# Remove the now-unused trailing rows 16-19 (reduces dimension from D19 to D15)
$ws.Rows("16:19").Delete()
